$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Add a new "Overdraft" column to the companies (Sheet2) table to support
# saving overdraft transaction data to file.
$ws2.Range("G1").Value = "Overdraft"
$ws2.Range("G3").Value = 0
$ws2.Range("G4").Value = 1500

# Reflect the selection state left behind in the saved workbook (entire
# row 5 selected on each sheet, via the row header).
$ws1.Select()
$ws1.Rows(5).Select()

$ws2.Select()
$ws2.Rows(5).Select()
